$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7

$ws.Cells.Item($row, 1).Value = 42612.894143518519
$ws.Cells.Item($row, 2).Value = 22
$ws.Cells.Item($row, 3).Value = 56
$ws.Cells.Item($row, 4).Value = 43
$ws.Cells.Item($row, 5).Value = 85
$ws.Cells.Item($row, 6).Value = 14
$ws.Cells.Item($row, 7).Value = 20934
$ws.Cells.Item($row, 8).Value = 15241
$ws.Cells.Item($row, 9).Value = 844
$ws.Cells.Item($row, 10).Value = 187
$ws.Cells.Item($row, 11).Value = 144
$ws.Cells.Item($row, 12).Value = 18
$ws.Cells.Item($row, 13).Value = 3
$ws.Cells.Item($row, 14).Value = "Named"
